$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.015.95'
$ws.Range('E2').Value = '  -1.52%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.569.84'
$ws.Range('E3').Value = '  -2.53%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.19'
$ws.Range('E5').Value = '  -3.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '186.08'
$ws.Range('E6').Value = '  -4.64%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.563.37'
$ws.Range('E7').Value = '  -2.55%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.621'
$ws.Range('E8').Value = '  -4.20%  '
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('E10').Value = '  -0.66%  '
$ws.Range('E11').Value = '  -3.99%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.97'
$ws.Range('E12').Value = '  -5.79%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000302'
$ws.Range('E13').Value = '  +1.74%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.53'
$ws.Range('E14').Value = '  -4.69%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.143.60'
$ws.Range('E15').Value = '  -2.44%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.64'
$ws.Range('E16').Value = '  -3.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.569.07'
$ws.Range('E17').Value = '  -2.61%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.902.88'
$ws.Range('E18').Value = '  -1.68%  '
$ws.Range('E19').Value = '  -1.82%  '
$ws.Range('E20').Value = '  -0.98%  '
$ws.Range('E21').Value = '  -3.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '490.82'
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.95'
$ws.Range('E23').Value = '  -1.30%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.92'
$ws.Range('E24').Value = '  -6.01%  '
$ws.Range('E25').Value = '  -2.98%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '95.34'
$ws.Range('E26').Value = '  +4.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.45'
$ws.Range('E27').Value = '  -0.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.96'
$ws.Range('E28').Value = '  -6.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.32'
$ws.Range('E29').Value = '  -3.31%  '
$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.63'
$ws.Range('E30').Value = '  -2.81%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '31.55'
$ws.Range('E31').Value = '  -4.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '66.83'
$ws.Range('E32').Value = '  +0.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.06'
$ws.Range('E33').Value = '  -1.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.114'
$ws.Range('E34').Value = '  -6.75%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '570.58'
$ws.Range('E35').Value = '  -8.45%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.26'
$ws.Range('E36').Value = '  +14.51%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '38.59'
$ws.Range('E37').Value = '  -4.53%  '
$ws.Range('B38').Value = 'TheGraph'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.403'
$ws.Range('E38').Value = '  -2.43%  '
$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  -0.09%  '
$ws.Range('B40').Value = 'PEPE'
$ws.Range('C40').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0790'
$ws.Range('E40').Value = '  -5.54%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.48'
$ws.Range('E41').Value = '  -3.25%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.15'
$ws.Range('E42').Value = '  -0.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.135'
$ws.Range('E43').Value = '  -9.57%  '
$ws.Range('E44').Value = '  -5.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.228.81'
$ws.Range('E45').Value = '  -3.22%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0442'
$ws.Range('E46').Value = '  -3.13%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.42'
$ws.Range('E47').Value = '  +2.67%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.57'
$ws.Range('E48').Value = '  -1.62%  '
$ws.Range('E49').Value = '  -2.59%  '
$ws.Range('B50').Value = 'OceanProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.52'
$ws.Range('E50').Value = '  +26.11%  '
$ws.Range('B51').Value = 'FirstDigitalUSD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.999'
$ws.Range('E51').Value = '  -0.07%  '
